$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the "twelve months ended" column headers by one fiscal year ---
# Old columns were 1396/12 .. 1400/12; new columns are 1397/12 .. 1401/12.
$ws.Range("E8").Value  = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value  = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value  = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value  = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value  = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Row 16: هزینه استهلاک ---
$ws.Range("E16").Value = 35086
$ws.Range("F16").Value = 39263
$ws.Range("G16").Value = 54080
$ws.Range("H16").Value = 100689
$ws.Range("I16").Value = 113045

# --- Row 17: هزینه حقوق و دستمزد ---
$ws.Range("E17").Value = 621730
$ws.Range("F17").Value = 948620
$ws.Range("G17").Value = 740888
$ws.Range("H17").Value = 1227299
$ws.Range("I17").Value = 6400127

# --- Row 19: سایر هزینه ها ---
$ws.Range("E19").Value = 4555459
$ws.Range("F19").Value = 6274903
$ws.Range("G19").Value = 7516679
$ws.Range("H19").Value = 11290910
$ws.Range("I19").Value = 19656535

# --- Row 20: جمع ---
$ws.Range("E20").Value = 5212275
$ws.Range("F20").Value = 7262786
$ws.Range("G20").Value = 8311647
$ws.Range("H20").Value = 12618898
$ws.Range("I20").Value = 26169707

# --- Row 26: تعداد پرسنل غیر تولیدی شرکت ---
$ws.Range("E26").Value = 630
$ws.Range("F26").Value = 630
$ws.Range("G26").Value = 575
$ws.Range("H26").Value = 690
$ws.Range("I26").Value = 1142

# --- Row 27: تعداد پرسنل تولیدی شرکت ---
$ws.Range("E27").Value = 402
$ws.Range("F27").Value = 446
$ws.Range("G27").Value = 504
$ws.Range("H27").Value = 517
$ws.Range("I27").Value = 2108
